# Add 2022-Q1 data: insert a new "2022-Q1" sheet before the "总计" (total)
# sheet, and refresh "总计" with a new summary row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper sheet used purely as a style donor: "2021-Q4" already has the
# 8-column per-fund layout with the bold/bordered header style (and the
# bold/bordered index column A) that the new sheet should reuse.
# ------------------------------------------------------------------
$styleDonor = $wb.Worksheets.Item("2021-Q4")

# The existing "总计" sheet - remember it so we can drop it once its
# replacement (with identical formatting, plus the new row) is ready.
$oldTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Delete()

# ------------------------------------------------------------------
# 1) Build the new "2022-Q1" per-fund holdings sheet.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($null, $lastSheet)
$q1.Name = "2022-Q1"

# Copy the header style (bold/centered/bordered) for B1:H1, and the
# index-column style (A2, bold/centered/bordered) down through row 14 -
# without ever touching A1, which the donor sheet leaves untouched too.
$styleDonor.Range("B1:H1").Copy($q1.Range("B1:H1"))
$styleDonor.Range("A2").Copy($q1.Range("A2:A14"))

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Header labels
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$q1Rows = @(
    @(0,  "004856", "广发中证全指建筑材料指数A", "13.72", "94.61", "3.39", "0.4651", 6),
    @(1,  "007994", "华夏中证500指数增强A", "31.45", "92.72", "1.40", "0.4403", 7),
    @(2,  "004857", "广发中证全指建筑材料指数C", "6.05", "94.61", "3.39", "0.2051", 6),
    @(3,  "003318", "景顺长城中证500行业中性低波动指数", "13.99", "93.88", "1.19", "0.1665", 7),
    @(4,  "159745", "国泰中证全指建筑材料交易型开放式指数证券投资基金", "3.76", "98.37", "3.54", "0.1331", 6),
    @(5,  "007995", "华夏中证500指数增强C", "5.45", "92.72", "1.40", "0.0763", 7),
    @(6,  "008856", "华夏安泰对冲策略3个月定期开放灵活配置混合", "2.87", "80.67", "1.31", "0.0376", 9),
    @(7,  "516750", "富国中证全指建筑材料ETF", "0.47", "98.22", "3.49", "0.0164", 7),
    @(8,  "512260", "华安中证500行业中性低波动ETF", "1.17", "96.94", "1.23", "0.0144", 7),
    @(9,  "164811", "工银瑞信中证京津冀协同发展主题指数（LOF）A", "0.23", "94.28", "2.97", "0.0068", 8),
    @(10, "512780", "广发中证京津冀协同发展主题ETF", "0.13", "98.52", "3.05", "0.0040", 9),
    @(11, "005126", "银河量化稳进混合", "0.10", "78.20", "2.05", "0.0020", 5),
    @(12, "164825", "工银瑞信中证京津冀协同发展主题指数（LOF）C", "0.06", "94.28", "2.97", "0.0018", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q1.Cells.Item($r, 2) $row[1]
    Set-TextValue $q1.Cells.Item($r, 3) $row[2]
    Set-TextValue $q1.Cells.Item($r, 4) $row[3]
    Set-TextValue $q1.Cells.Item($r, 5) $row[4]
    Set-TextValue $q1.Cells.Item($r, 6) $row[5]
    Set-TextValue $q1.Cells.Item($r, 7) $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$q1.Range("A1").Select()

# ------------------------------------------------------------------
# 2) Rebuild the "总计" (grand total) sheet with the 2022-Q1 row added
#    on top and every other row shifted down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$styleDonor.Range("B1:D1").Copy($total.Range("B1:D1"))
$styleDonor.Range("A2").Copy($total.Range("A2:A7"))

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 13, 1.57),
    @(1, "2021-Q4", 7, 0.92),
    @(2, "2021-Q3", 6, 0.46),
    @(3, "2021-Q2", 4, 0.49),
    @(4, "2021-Q1", 5, 0.52),
    @(5, "2020-Q4", 2, 0.67)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$total.Range("A1").Select()
